$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume columns in this sheet are plain text (inline strings) in the
# source workbook, even when the text happens to look like a plain number (e.g.
# "242.94"), because many of them use "." as a thousands separator (e.g.
# "29.327.37") and must stay textual and consistent column-to-column. Forcing the
# NumberFormat to Text before assigning the value keeps Excel from auto-converting
# those into real numbers; resetting the Style back to Normal afterwards avoids
# leaving a stray number-format style on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.327.37"
$ws.Range("E2").Value = "  -0.01%  "
Set-TextValue $ws.Range("D3") "1.876.83"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.05%  "
Set-TextValue $ws.Range("D6") "242.94"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue $ws.Range("D8") "0.08004"
$ws.Range("E8").Value = "  +2.92%  "
Set-TextValue $ws.Range("D9") "0.3156"
$ws.Range("E9").Value = "  +1.55%  "
Set-TextValue $ws.Range("D10") "25.01"
$ws.Range("E10").Value = "  -0.36%  "
Set-TextValue $ws.Range("D11") "0.08244"
$ws.Range("E11").Value = "  -1.85%  "
Set-TextValue $ws.Range("D12") "1.896.35"
$ws.Range("E12").Value = "  +1.30%  "
Set-TextValue $ws.Range("D13") "5.249"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  +3.91%  "
Set-TextValue $ws.Range("D15") "0.7128"
$ws.Range("E15").Value = "  +0.30%  "
Set-TextValue $ws.Range("D16") "6.392"
$ws.Range("E16").Value = "  +5.22%  "
Set-TextValue $ws.Range("D17") "0.000008565"
$ws.Range("E17").Value = "  +4.64%  "
Set-TextValue $ws.Range("D18") "29.357.37"
$ws.Range("E18").Value = "  +0.08%  "
Set-TextValue $ws.Range("D19") "244.20"
$ws.Range("E19").Value = "  +1.73%  "
Set-TextValue $ws.Range("D20") "2.155.46"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("E22").Value = "  +0.02%  "
Set-TextValue $ws.Range("D23") "7.783"
$ws.Range("E23").Value = "  +0.34%  "
Set-TextValue $ws.Range("D24") "1.002"
$ws.Range("E24").Value = "  +0.11%  "
Set-TextValue $ws.Range("D25") "0.1561"
$ws.Range("E25").Value = "  -2.42%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D26") "9.049"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D27") "162.53"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("E30").Value = "  +0.32%  "
Set-TextValue $ws.Range("D31") "4.315"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  -7.70%  "
Set-TextValue $ws.Range("D33") "0.05380"
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("E36").Value = "  +0.47%  "
Set-TextValue $ws.Range("D37") "2.690"
$ws.Range("E37").Value = "  -0.43%  "
Set-TextValue $ws.Range("D38") "0.01879"
$ws.Range("E38").Value = "  +0.42%  "
Set-TextValue $ws.Range("D39") "1.253.89"
$ws.Range("E39").Value = "  +2.72%  "
Set-TextValue $ws.Range("D40") "2.752"
$ws.Range("E40").Value = "  +1.03%  "
Set-TextValue $ws.Range("D41") "6.490"
$ws.Range("E41").Value = "  -0.64%  "
Set-TextValue $ws.Range("D42") "0.9183"
$ws.Range("E42").Value = "  +3.53%  "
Set-TextValue $ws.Range("D43") "112.77"
$ws.Range("E43").Value = "  +2.65%  "
Set-TextValue $ws.Range("D44") "74.21"
$ws.Range("E44").Value = "  +2.49%  "
Set-TextValue $ws.Range("D45") "0.00000000134"
$ws.Range("E45").Value = "  +9.10%  "
$ws.Range("E46").Value = "  +0.00%  "
Set-TextValue $ws.Range("D47") "2.045.31"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  +0.23%  "
Set-TextValue $ws.Range("D50") "9.475"
$ws.Range("E50").Value = "  +1.09%  "
Set-TextValue $ws.Range("D51") "0.4360"
$ws.Range("E51").Value = "  +1.10%  "
